# Daily attendance processing - 2025-10-18 13:04:28
# Normalize the "Recorded By" (column G) ordering on the Session Analysis
# Results sheet: the author swapped ordering between the "System" /
# "system" recorder token and the human/service-account email token for a
# batch of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
    elseif ($val -eq "System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, System"
    }
    elseif ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
